$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells are plain text (inline strings) holding values such
# as "298.87" or "0.99%". Setting .Value directly with such strings
# makes Excel auto-detect them as numbers/percentages, which would
# change the cell type. Forcing NumberFormat "@" (Text) on each cell
# immediately before assigning its value keeps the cell as Text,
# matching the original data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "298.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.99%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.73%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.142"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.41%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07949"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "7.93%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.541"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "54.31%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.81%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.826"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.19%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9071"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.31%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1738"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.86%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07243"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.74%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08052"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.37%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03021"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09962"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.56%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001504"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.76%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005960"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.02%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.502"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.58%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.254"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.18%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1318"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.01%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.627"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.22%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.26%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04578"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.80%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001258"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.41%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004455"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.73%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-9.12%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003429"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "83.12%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01851"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "10.12%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04522"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.45%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007033"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.95%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.30%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002239"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.81%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01044"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.35%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006453"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.70%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.06%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-57.44%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.006198"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-39.32%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
